$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CodeSchemes")

# Insert a new column before column B, shifting B:W to C:X
$ws.Columns("B:B").Insert()

# Set the new column's header and the data row value
$ws.Range("B1").Value = "ORGANIZATION"
$ws.Range("B2").Value = "74a41211-8c99-4835-a519-7a61612b1098"
